$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 96 (shifting the existing rows 96.. down to 98..)
$ws.Rows.Item(96).EntireRow.Insert()
$ws.Rows.Item(96).EntireRow.Insert()

# New row 96: Femacal de La Calera / Papaya / Cultivar IV Región / Primera
$ws.Range("A96").Value = 3
$ws.Range("B96").Value = "Femacal de La Calera"
$ws.Range("C96").Value = "Coquimbo"
$ws.Range("D96").Value = 45211
$ws.Range("E96").Value = 5
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100108
$ws.Range("H96").Value = "Tropicales y subtropicales"
$ws.Range("I96").Value = 100108004
$ws.Range("J96").Value = "Papaya"
$ws.Range("K96").Value = "Cultivar IV Región"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 60
$ws.Range("N96").Value = 16000
$ws.Range("O96").Value = 16000
$ws.Range("P96").Value = 16000
$ws.Range("Q96").Value = "$/bandeja 10 kilos"
$ws.Range("R96").Value = "Provincia del Elquí"
$ws.Range("S96").Value = 1600
$ws.Range("T96").Value = 10

# New row 97: Femacal de La Calera / Papaya / Cultivar IV Región / Segunda
$ws.Range("A97").Value = 3
$ws.Range("B97").Value = "Femacal de La Calera"
$ws.Range("C97").Value = "Coquimbo"
$ws.Range("D97").Value = 45211
$ws.Range("E97").Value = 5
$ws.Range("F97").Value = "Fruta"
$ws.Range("G97").Value = 100108
$ws.Range("H97").Value = "Tropicales y subtropicales"
$ws.Range("I97").Value = 100108004
$ws.Range("J97").Value = "Papaya"
$ws.Range("K97").Value = "Cultivar IV Región"
$ws.Range("L97").Value = "Segunda"
$ws.Range("M97").Value = 57
$ws.Range("N97").Value = 13000
$ws.Range("O97").Value = 13000
$ws.Range("P97").Value = 13000
$ws.Range("Q97").Value = "$/bandeja 10 kilos"
$ws.Range("R97").Value = "Provincia del Elquí"
$ws.Range("S97").Value = 1300
$ws.Range("T97").Value = 10
